# V 2.5.9 - Work on Wood-Production recipe table
# Increase the "Time" inputs (column D) for Water, Ash and Fertilizer rows,
# which cascades through the derived formulas in columns G/J/L.
# Also simplify the Water row's Wood "Time" formula (J4) and flatten the
# now-static Ash/Fertilizer "Time" values in column J, leaving the
# selection on J7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the base "Time" (D) inputs that drive the recipe table
$ws.Range("D4").Value = 200
$ws.Range("D5").Value = 150
$ws.Range("D6").Value = 100

# Water row: Wood Time is now simply equal to the Seedlings Time
$ws.Range("J4").Formula = "=G4"

# Ash / Fertilizer rows: Wood Time becomes a static value instead of a formula
$ws.Range("J5").Value = 360
$ws.Range("J6").Value = 300

# Leave the active selection on J7, matching the author's last edit location
$ws.Range("J7").Select()
